$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# "adicionando usuario ao banco": update the status of requisito 2 (row 3)
# from "Pendente" (shared string 21) to "Pronto" (shared string 20).
$ws.Range("C3").Value = "Pronto"

# Update the active cell / selection to C4, matching the saved cursor position.
[void]$ws.Range("C4").Select()
